$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the A3/B3 registration values and switch C3 to a numeric literal
$ws.Range("A3").Value = "avu"
$ws.Range("B3").Value = "Устюжанин Александр Викторович"
$ws.Range("C3").Value = 1

# Update the current selection shown in the saved view
$ws.Range("C2").Select()
